$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 16.46174698877879
$ws.Range("C2").Value = 7.56256815758016
$ws.Range("D2").Value = 6.922430688434947
$ws.Range("F2").Value = 46.113983916542
$ws.Range("G2").Value = 56.62322186534746
$ws.Range("H2").Value = 21.23646667299198
$ws.Range("J2").Value = 11.09443381621367
$ws.Range("K2").Value = 12.40524912010419
$ws.Range("L2").Value = 11.36053880536494
$ws.Range("B3").Value = 16.33619126908918
$ws.Range("C3").Value = 7.528117260735494
$ws.Range("D3").Value = 6.923021228606737
$ws.Range("F3").Value = 46.11006445916851
$ws.Range("G3").Value = 56.55961597645936
$ws.Range("H3").Value = 21.2637083512654
$ws.Range("J3").Value = 11.1140213616196
$ws.Range("K3").Value = 12.32273277723724
$ws.Range("L3").Value = 11.36869780997093
$ws.Range("B4").Value = 16.26266156905998
$ws.Range("C4").Value = 7.506309983433146
$ws.Range("D4").Value = 6.924382918595998
$ws.Range("F4").Value = 46.11585245550614
$ws.Range("G4").Value = 56.5318498835216
$ws.Range("H4").Value = 21.28337192444215
$ws.Range("J4").Value = 11.12696104233393
$ws.Range("K4").Value = 12.27474126608153
$ws.Range("L4").Value = 11.37533001162381
$ws.Range("B5").Value = 16.23362360711646
$ws.Range("C5").Value = 7.49725983194435
$ws.Range("D5").Value = 6.925190022033284
$ws.Range("F5").Value = 46.12027240367318
$ws.Range("G5").Value = 56.52337799075548
$ws.Range("H5").Value = 21.29212311361967
$ws.Range("J5").Value = 11.13246400809068
$ws.Range("K5").Value = 12.25587559783505
$ws.Range("L5").Value = 11.37844134990504
$ws.Range("B6").Value = 16.22885865552789
$ws.Range("C6").Value = 7.495747215539879
$ws.Range("D6").Value = 6.92533930042119
$ws.Range("F6").Value = 46.1211308321632
$ws.Range("G6").Value = 56.52214304336728
$ws.Range("H6").Value = 21.29362080824618
$ws.Range("J6").Value = 11.13339167091281
$ws.Range("K6").Value = 12.25278524840441
$ws.Range("L6").Value = 11.37898268509551
$ws.Range("B7").Value = 16.26226616488623
$ws.Range("C7").Value = 7.506188590610738
$ws.Range("D7").Value = 6.924392781071659
$ws.Range("F7").Value = 46.11590371820606
$ws.Range("G7").Value = 56.53172411297117
$ws.Range("H7").Value = 21.28348695808412
$ws.Range("J7").Value = 11.12703432570588
$ws.Range("K7").Value = 12.27448401407493
$ws.Range("L7").Value = 11.37537031674994
$ws.Range("B8").Value = 16.41773474678078
$ws.Range("C8").Value = 7.550825887785667
$ws.Range("D8").Value = 6.922427548520198
$ws.Range("F8").Value = 46.11093238747118
$ws.Range("G8").Value = 56.59895097121743
$ws.Range("H8").Value = 21.2452497841993
$ws.Range("J8").Value = 11.10099837925801
$ws.Range("K8").Value = 12.3762542052473
$ws.Range("L8").Value = 11.36301569221146
$ws.Range("B9").Value = 16.74946504207943
$ws.Range("C9").Value = 7.633125890815204
$ws.Range("D9").Value = 6.926455987302252
$ws.Range("F9").Value = 46.16614060537584
$ws.Range("G9").Value = 56.82007255011887
$ws.Range("H9").Value = 21.193588051249
$ws.Range("J9").Value = 11.0571673045664
$ws.Range("K9").Value = 12.59615891188451
$ws.Range("L9").Value = 11.35163099881998
$ws.Range("B10").Value = 17.00758303488314
$ws.Range("C10").Value = 7.690348463876758
$ws.Range("D10").Value = 6.93415683102782
$ws.Range("F10").Value = 46.24616904561855
$ws.Range("G10").Value = 57.03647843685992
$ws.Range("H10").Value = 21.16986919021863
$ws.Range("J10").Value = 11.02934520948678
$ws.Range("K10").Value = 12.76887081832503
$ws.Range("L10").Value = 11.3510501139327
$ws.Range("B11").Value = 17.12769100841578
$ws.Range("C11").Value = 7.715665419824752
$ws.Range("D11").Value = 6.938674438912035
$ws.Range("F11").Value = 46.29109789072588
$ws.Range("G11").Value = 57.14649599244773
$ws.Range("H11").Value = 21.16217238834615
$ws.Range("J11").Value = 11.01763432141696
$ws.Range("K11").Value = 12.84958043644392
$ws.Range("L11").Value = 11.35246436605968
$ws.Range("B12").Value = 17.17352122973022
$ws.Range("C12").Value = 7.725148657150307
$ws.Range("D12").Value = 6.940529587412295
$ws.Range("F12").Value = 46.30933100482898
$ws.Range("G12").Value = 57.18980515703432
$ws.Range("H12").Value = 21.15970250236079
$ws.Range("J12").Value = 11.01333528250484
$ws.Range("K12").Value = 12.88042650867064
$ws.Range("L12").Value = 11.35324010437206
$ws.Range("B13").Value = 17.16363601732881
$ws.Range("C13").Value = 7.723110903570713
$ws.Range("D13").Value = 6.940123648417967
$ws.Range("F13").Value = 46.30535004682815
$ws.Range("G13").Value = 57.18040474423219
$ws.Range("H13").Value = 21.16021465932196
$ws.Range("J13").Value = 11.01425513117916
$ws.Range("K13").Value = 12.87377106590453
$ws.Range("L13").Value = 11.35306237243349
$ws.Range("B14").Value = 17.13145473821186
$ws.Range("C14").Value = 7.716447696848507
$ws.Range("D14").Value = 6.938824177117018
$ws.Range("F14").Value = 46.29257353745822
$ws.Range("G14").Value = 57.1500261378232
$ws.Range("H14").Value = 21.16196027710503
$ws.Range("J14").Value = 11.01727792083549
$ws.Range("K14").Value = 12.85211263208903
$ws.Range("L14").Value = 11.35252337980404
$ws.Range("B15").Value = 17.11178693172
$ws.Range("C15").Value = 7.712352747661158
$ws.Range("D15").Value = 6.938046977907344
$ws.Range("F15").Value = 46.28490617657785
$ws.Range("G15").Value = 57.13163245044019
$ws.Range("H15").Value = 21.1630874312132
$ws.Range("J15").Value = 11.01914711833702
$ws.Range("K15").Value = 12.8388823356633
$ws.Range("L15").Value = 11.35222447472643
$ws.Range("B16").Value = 16.99978450297815
$ws.Range("C16").Value = 7.688679482731884
$ws.Range("D16").Value = 6.93388188612285
$ws.Range("F16").Value = 46.24340385364151
$ws.Range("G16").Value = 57.02952006735156
$ws.Range("H16").Value = 21.17043443279079
$ws.Range("J16").Value = 11.03012954038472
$ws.Range("K16").Value = 12.76363723743088
$ws.Range("L16").Value = 11.35099136684501
$ws.Range("B17").Value = 16.93173591332103
$ws.Range("C17").Value = 7.673973080185878
$ws.Range("D17").Value = 6.931585557990429
$ws.Range("F17").Value = 46.22012255769442
$ws.Range("G17").Value = 56.96983068085055
$ws.Range("H17").Value = 21.17573378556052
$ws.Range("J17").Value = 11.03710882575043
$ws.Range("K17").Value = 12.71800784795778
$ws.Range("L17").Value = 11.35066399555579
$ws.Range("B18").Value = 16.89285159899267
$ws.Range("C18").Value = 7.665447217009731
$ws.Range("D18").Value = 6.930360347938668
$ws.Range("F18").Value = 46.20753447667496
$ws.Range("G18").Value = 56.93658937861409
$ws.Range("H18").Value = 21.17907297425825
$ws.Range("J18").Value = 11.04121214364716
$ws.Range("K18").Value = 12.69196617961607
$ws.Range("L18").Value = 11.35063376523717
$ws.Range("B19").Value = 16.87973108152103
$ws.Range("C19").Value = 7.662549010111388
$ws.Range("D19").Value = 6.929961971169391
$ws.Range("F19").Value = 46.2034104042604
$ws.Range("G19").Value = 56.92552219610972
$ws.Range("H19").Value = 21.18025356934134
$ws.Range("J19").Value = 11.04261675596827
$ws.Range("K19").Value = 12.68318456207243
$ws.Range("L19").Value = 11.35065071171447
$ws.Range("B20").Value = 16.93895364436861
$ws.Range("C20").Value = 7.675545557427937
$ws.Range("D20").Value = 6.931820124415148
$ws.Range("F20").Value = 46.222517855599
$ws.Range("G20").Value = 56.9760719718775
$ws.Range("H20").Value = 21.175139528297
$ws.Range("J20").Value = 11.03635665842734
$ws.Range("K20").Value = 12.72284432844691
$ws.Range("L20").Value = 11.3506824938441
$ws.Range("B21").Value = 17.14089802745399
$ws.Range("C21").Value = 7.718407666386205
$ws.Range("D21").Value = 6.939201955534431
$ws.Range("F21").Value = 46.29629326407009
$ws.Range("G21").Value = 57.15890448752096
$ws.Range("H21").Value = 21.16143547813351
$ws.Range("J21").Value = 11.01638637628447
$ws.Range("K21").Value = 12.85846675685418
$ws.Range("L21").Value = 11.35267518613163
$ws.Range("B22").Value = 17.27489022264071
$ws.Range("C22").Value = 7.745815335459932
$ws.Range("D22").Value = 6.944867520242489
$ws.Range("F22").Value = 46.35161457960321
$ws.Range("G22").Value = 57.28799235460696
$ws.Range("H22").Value = 21.15507117806309
$ws.Range("J22").Value = 11.00412498224423
$ws.Range("K22").Value = 12.94874210285156
$ws.Range("L22").Value = 11.35537693066937
$ws.Range("B23").Value = 17.20320504205378
$ws.Range("C23").Value = 7.731243078319733
$ws.Range("D23").Value = 6.941767239190963
$ws.Range("F23").Value = 46.32144075099253
$ws.Range("G23").Value = 57.21822359175511
$ws.Range("H23").Value = 21.15823079771623
$ws.Range("J23").Value = 11.01059691632402
$ws.Range("K23").Value = 12.90041888921451
$ws.Range("L23").Value = 11.35380732725588
$ws.Range("B24").Value = 16.93568976471537
$ws.Range("C24").Value = 7.674834861421067
$ws.Range("D24").Value = 6.931713780989758
$ws.Range("F24").Value = 46.22143245936492
$ws.Range("G24").Value = 56.97324693455504
$ws.Range("H24").Value = 21.17540728076159
$ws.Range("J24").Value = 11.03669643030517
$ws.Range("K24").Value = 12.7206571579091
$ws.Range("L24").Value = 11.35067363859227
$ws.Range("B25").Value = 16.65704755311648
$ws.Range("C25").Value = 7.61142732957915
$ws.Range("D25").Value = 6.924528639797403
$ws.Range("F25").Value = 46.1442609837455
$ws.Range("G25").Value = 56.75073616210804
$ws.Range("H25").Value = 21.20506471903314
$ws.Range("J25").Value = 11.06825375442079
$ws.Range("K25").Value = 12.53461747177074
$ws.Range("L25").Value = 11.35334031065788
